$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 271.72726
$ws.Range("I2").Value = 280.9
$ws.Range("K2").Value = 280.9
$ws.Range("M2").Value = -167.9

$ws.Range("H86").Value = 4099.9
$ws.Range("I86").Value = 5340
$ws.Range("J86").Value = 2859.8
$ws.Range("K86").Value = 5340
$ws.Range("L86").Value = 2859.8
$ws.Range("M86").Value = -4217
$ws.Range("N86").Value = -5105.8

$ws.Range("H89").Value = 4099.9
$ws.Range("I89").Value = 5340
$ws.Range("J89").Value = 2859.8
$ws.Range("K89").Value = 26700
$ws.Range("L89").Value = 14299
$ws.Range("M89").Value = -21084
$ws.Range("N89").Value = -25531

$ws.Range("H96").Value = 1087.0526
$ws.Range("I96").Value = 1529.8182
$ws.Range("J96").Value = 478.25
$ws.Range("K96").Value = 4589.4546
$ws.Range("L96").Value = 1434.75
$ws.Range("M96").Value = -3216.4546
$ws.Range("N96").Value = -4180.75

$ws.Range("H100").Value = 1110.8
$ws.Range("I100").Value = 766.5
$ws.Range("J100").Value = 1914.1666
$ws.Range("K100").Value = 766.5
$ws.Range("L100").Value = 1914.1666
$ws.Range("M100").Value = -225.5
$ws.Range("N100").Value = -2996.1666

$ws.Range("H112").Value = 1914.421
$ws.Range("I112").Value = 729.6667
$ws.Range("J112").Value = 2136.5625
$ws.Range("K112").Value = 2189.0001
$ws.Range("L112").Value = 6409.6875
$ws.Range("M112").Value = -1081.0001
$ws.Range("N112").Value = -8625.6875

$ws.Range("H138").Value = 1960.6771
$ws.Range("I138").Value = 788.2857
$ws.Range("J138").Value = 2160.8416
$ws.Range("K138").Value = 2364.8571
$ws.Range("L138").Value = 6482.524800000001
$ws.Range("M138").Value = 2775.1429
$ws.Range("N138").Value = -16762.5248

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 993
$ws.Range("I2").Value = 708.2353000000001
$ws.Range("K2").Value = 708.2353000000001
$ws.Range("M2").Value = -595.2353000000001

$ws.Range("H4").Value = 707.5
$ws.Range("I4").Value = 499
$ws.Range("J4").Value = 916
$ws.Range("K4").Value = 499
$ws.Range("L4").Value = 916
$ws.Range("M4").Value = -383
$ws.Range("N4").Value = -1148

$ws.Range("H110").Value = 426.2
$ws.Range("I110").Value = 426.2
$ws.Range("K110").Value = 426.2
$ws.Range("M110").Value = 1618.8

$ws.Range("H116").Value = 993
$ws.Range("I116").Value = 708.2353000000001
$ws.Range("K116").Value = 708.2353000000001
$ws.Range("M116").Value = 1585.7647

$ws.Range("H122").Value = 2627.0386
$ws.Range("I122").Value = 2561.5417
$ws.Range("J122").Value = 3413
$ws.Range("K122").Value = 7684.625100000001
$ws.Range("L122").Value = 10239
$ws.Range("M122").Value = -5234.625100000001
$ws.Range("N122").Value = -15139

$ws.Range("H132").Value = 3019.0513
$ws.Range("I132").Value = 2207.56
$ws.Range("K132").Value = 6622.68
$ws.Range("M132").Value = -4092.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 993
$ws.Range("I3").Value = 708.2353000000001
$ws.Range("K3").Value = 708.2353000000001
$ws.Range("M3").Value = -594.2353000000001

$ws.Range("H20").Value = 2639
$ws.Range("I20").Value = 2833
$ws.Range("J20").Value = 2289.8
$ws.Range("K20").Value = 2833
$ws.Range("L20").Value = 2289.8
$ws.Range("M20").Value = -2586
$ws.Range("N20").Value = -2783.8

$ws.Range("H134").Value = 9032.538
$ws.Range("I134").Value = 1451.9166
$ws.Range("J134").Value = 100000
$ws.Range("K134").Value = 4355.7498
$ws.Range("L134").Value = 300000
$ws.Range("M134").Value = -1820.7498
$ws.Range("N134").Value = -305070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1711.3903
$ws.Range("I31").Value = 1711.3903
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1711.3903
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1416.3903
$ws.Range("N31").ClearContents()

$ws.Range("H34").Value = 1711.3903
$ws.Range("I34").Value = 1711.3903
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1711.3903
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1509.3903
$ws.Range("N34").ClearContents()

$ws.Range("H86").Value = 2922591.8
$ws.Range("I86").Value = 4459993.5
$ws.Range("J86").Value = 39963.5
$ws.Range("K86").Value = 4459993.5
$ws.Range("L86").Value = 39963.5
$ws.Range("M86").Value = -4458870.5
$ws.Range("N86").Value = -42209.5

$ws.Range("H89").Value = 2922591.8
$ws.Range("I89").Value = 4459993.5
$ws.Range("J89").Value = 39963.5
$ws.Range("K89").Value = 22299967.5
$ws.Range("L89").Value = 199817.5
$ws.Range("M89").Value = -22294351.5
$ws.Range("N89").Value = -211049.5

$ws.Range("H99").Value = 1893.875
$ws.Range("I99").Value = 1741.8334
$ws.Range("J99").Value = 2350
$ws.Range("K99").Value = 1741.8334
$ws.Range("L99").Value = 2350
$ws.Range("M99").Value = -243.8334
$ws.Range("N99").Value = -5346

$ws.Range("H126").Value = 1893.875
$ws.Range("I126").Value = 1741.8334
$ws.Range("J126").Value = 2350
$ws.Range("K126").Value = 5225.5002
$ws.Range("L126").Value = 7050
$ws.Range("M126").Value = -2755.5002
$ws.Range("N126").Value = -11990

$ws.Range("H132").Value = 2320.5454
$ws.Range("I132").Value = 1796.7142
$ws.Range("J132").Value = 3237.25
$ws.Range("K132").Value = 5390.142599999999
$ws.Range("L132").Value = 9711.75
$ws.Range("M132").Value = -2860.142599999999
$ws.Range("N132").Value = -14771.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 14325.363

$ws.Range("H66").Value = 14325.363

$ws.Range("H87").Value = 1598.1428
$ws.Range("J87").Value = 1899.8
$ws.Range("L87").Value = 5699.4
$ws.Range("N87").Value = -8195.4

$ws.Range("H90").Value = 1598.1428
$ws.Range("J90").Value = 1899.8
$ws.Range("L90").Value = 17098.2
$ws.Range("N90").Value = -29578.2

$ws.Range("H113").Value = 759.2174
$ws.Range("I113").Value = 723
$ws.Range("J113").Value = 760.86365
$ws.Range("K113").Value = 2169
$ws.Range("L113").Value = 2282.59095
$ws.Range("M113").Value = 1
$ws.Range("N113").Value = -6622.59095

$ws.Range("H131").Value = 16975416
$ws.Range("J131").Value = 34274.043
$ws.Range("L131").Value = 102822.129
$ws.Range("N131").Value = -112902.129

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1008.8333
$ws.Range("I97").Value = 822.8889
$ws.Range("J97").Value = 1566.6666
$ws.Range("K97").Value = 822.8889
$ws.Range("L97").Value = 1566.6666
$ws.Range("M97").Value = -326.8889
$ws.Range("N97").Value = -2558.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1772.8667
$ws.Range("I16").Value = 1185.2142
$ws.Range("K16").Value = 1185.2142
$ws.Range("M16").Value = -1015.2142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1815.8837
$ws.Range("I132").Value = 1402.5714
$ws.Range("J132").Value = 3624.125
$ws.Range("K132").Value = 4207.7142
$ws.Range("L132").Value = 10872.375
$ws.Range("M132").Value = -1677.7142
$ws.Range("N132").Value = -15932.375
